$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Move the two footer rows (old 24/25 -> new 29/30), keeping their
#    formatting (font/border/merge) intact, then clear the old locations.
# ---------------------------------------------------------------------------
$ws.Range("B24:C24").Copy($ws.Range("B29:C29"))
$ws.Range("H24:J24").Copy($ws.Range("H29:J29"))
$ws.Range("B25:C25").Copy($ws.Range("B30:C30"))
$ws.Range("H25:J25").Copy($ws.Range("H30:J30"))

$ws.Range("B24:C24").UnMerge()
$ws.Range("H24:J24").UnMerge()
$ws.Range("B25:C25").UnMerge()
$ws.Range("H25:J25").UnMerge()
$ws.Range("B24:J25").Clear()

# ---------------------------------------------------------------------------
# 2) Stamp out the new data-row formatting template (rows 16-23 use the
#    plain "middle of table" style; row 24 needs the heavier bottom-border
#    "end of table" style that row 19 currently carries). Grab the
#    end-of-table style from row 19 into row 24 FIRST, then re-stamp row 19
#    itself with the plain middle-of-table style copied from row 18.
# ---------------------------------------------------------------------------
$ws.Range("B19:J19").Copy($ws.Range("B24:J24"))
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))
$ws.Range("B18:J18").Copy($ws.Range("B20:J20"))
$ws.Range("B18:J18").Copy($ws.Range("B21:J21"))
$ws.Range("B18:J18").Copy($ws.Range("B22:J22"))
$ws.Range("B18:J18").Copy($ws.Range("B23:J23"))

# ---------------------------------------------------------------------------
# 3) Header / summary cells (labels unchanged, totals updated)
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 529308
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("C13").Value = 4
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 3
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# ---------------------------------------------------------------------------
# 4) Data rows 16-24 (new account-statement detail)
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047459638"
$ws.Range("D16").Value = "EDILSON ALBERTO CUERVO CASTRILLON"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047459638"
$ws.Range("D17").Value = "EDILSON ALBERTO CUERVO CASTRILLON"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047459638"
$ws.Range("D18").Value = "EDILSON ALBERTO CUERVO CASTRILLON"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "7960265"
$ws.Range("D19").Value = "WILSON ENRIQUE CASTILLA RODRIGUEZ"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 64940
$ws.Range("G19").Value = 1623500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "7960265"
$ws.Range("D20").Value = "WILSON ENRIQUE CASTILLA RODRIGUEZ"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 64940
$ws.Range("G20").Value = 1623500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "7960265"
$ws.Range("D21").Value = "WILSON ENRIQUE CASTILLA RODRIGUEZ"
$ws.Range("E21").Value = "2505"
$ws.Range("F21").Value = 64940
$ws.Range("G21").Value = 1623500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047477963"
$ws.Range("D22").Value = "YONEL SALCEDO PEREZ"
$ws.Range("E22").Value = "2507"
$ws.Range("F22").Value = 49788
$ws.Range("G22").Value = 908526

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1045231371"
$ws.Range("D23").Value = "OMAR LUIS LOPEZ HERNANDEZ"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1045231371"
$ws.Range("D24").Value = "OMAR LUIS LOPEZ HERNANDEZ"
$ws.Range("E24").Value = "2506"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500
